# BDD_materiaux.xlsx — rename the "stretched pipes" material label to be more
# descriptive, then widen column A so the longer text is fully visible, and
# leave the selection on A6 (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Conduites étirées" -> "Conduites étirées : aluminium, cuivre, …"
$ws.Range("A2").Value = "Conduites étirées : aluminium, cuivre, …"

# Column A now holds longer text: autosize it to fit the new content.
$ws.Columns.Item(1).AutoFit()

# Leave the active selection on A6 (single cell), matching the saved view.
$ws.Range("A6").Select() | Out-Null

# Nudge the saved window position (xWindow 230 -> 340).
$excel.Left = 340
